$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 9.70568990943505
$ws.Range("C2").Value = 5.38013134452062
$ws.Range("D2").Value = 5.98413062061643
$ws.Range("E2").Value = 16.43976262094679
$ws.Range("G2").Value = 3.634250864669381
$ws.Range("I2").Value = 20.91557936750216
$ws.Range("K2").Value = 9.724716203966702
$ws.Range("N2").Value = 18.18580905543253
$ws.Range("O2").Value = 22.68909384462043
# Row 3
$ws.Range("B3").Value = 9.383033473804028
$ws.Range("C3").Value = 5.058141228933312
$ws.Range("D3").Value = 5.865213763300594
$ws.Range("E3").Value = 15.51140536638115
$ws.Range("G3").Value = 3.636354570984811
$ws.Range("I3").Value = 20.98411211858861
$ws.Range("K3").Value = 9.496281678414075
$ws.Range("N3").Value = 18.24924469246868
$ws.Range("O3").Value = 22.72587309930514
# Row 4
$ws.Range("B4").Value = 9.181108907668078
$ws.Range("C4").Value = 4.848741193276935
$ws.Range("D4").Value = 5.792747955805352
$ws.Range("E4").Value = 14.91694087007511
$ws.Range("G4").Value = 3.63771358871603
$ws.Range("I4").Value = 21.03058583451466
$ws.Range("K4").Value = 9.35521212778198
$ws.Range("N4").Value = 18.28992947898172
$ws.Range("O4").Value = 22.75399825637059
# Row 5
$ws.Range("B5").Value = 9.097999931173433
$ws.Range("C5").Value = 4.760480703275369
$ws.Range("D5").Value = 5.763401468553383
$ws.Range("E5").Value = 14.66881527202015
$ws.Range("G5").Value = 3.63828438867105
$ws.Range("I5").Value = 21.05062628659826
$ws.Range("K5").Value = 9.297608629784202
$ws.Range("N5").Value = 18.30694659239406
$ws.Range("O5").Value = 22.76684898986149
# Row 6
$ws.Range("B6").Value = 9.084154211303117
$ws.Range("C6").Value = 4.745648727595729
$ws.Range("D6").Value = 5.758540958835388
$ws.Range("E6").Value = 14.62726825883914
$ws.Range("G6").Value = 3.638380197368365
$ws.Range("I6").Value = 21.05402046112233
$ws.Range("K6").Value = 9.288039182176412
$ws.Range("N6").Value = 18.30979875138987
$ws.Range("O6").Value = 22.76906663483816
# Row 7
$ws.Range("B7").Value = 9.179991216465009
$ws.Range("C7").Value = 4.847562717388568
$ws.Range("D7").Value = 5.792351373415143
$ws.Range("E7").Value = 14.91361795754272
$ws.Range("G7").Value = 3.637721217859776
$ws.Range("I7").Value = 21.03085164895879
$ws.Range("K7").Value = 9.354435621236004
$ws.Range("N7").Value = 18.2901572034773
$ws.Range("O7").Value = 22.75416594578405
# Row 8
$ws.Range("B8").Value = 9.595307115374315
$ws.Range("C8").Value = 5.271543221790744
$ws.Range("D8").Value = 5.943041273664996
$ws.Range("E8").Value = 16.12488676305872
$ws.Range("G8").Value = 3.634962282280517
$ws.Range("I8").Value = 20.93829565810788
$ws.Range("K8").Value = 9.646171327271023
$ws.Range("N8").Value = 18.20732242559192
$ws.Range("O8").Value = 22.70062257398645
# Row 9
$ws.Range("B9").Value = 10.37409826290512
$ws.Range("C9").Value = 6.009804477846066
$ws.Range("D9").Value = 6.240930904269886
$ws.Range("E9").Value = 18.36371166574409
$ws.Range("G9").Value = 3.630083646805282
$ws.Range("I9").Value = 20.79179092441747
$ws.Range("K9").Value = 10.20818637094255
$ws.Range("N9").Value = 18.05858346619118
$ws.Range("O9").Value = 22.63977201108549
# Row 10
$ws.Range("B10").Value = 10.91810892605774
$ws.Range("C10").Value = 6.495176272043225
$ws.Range("D10").Value = 6.458681962629338
$ws.Range("E10").Value = 19.99592572884205
$ws.Range("G10").Value = 3.626819718980564
$ws.Range("I10").Value = 20.70566452177062
$ws.Range("K10").Value = 10.6103514687513
$ws.Range("N10").Value = 17.9575609774049
$ws.Range("O10").Value = 22.62218413524824
# Row 11
$ws.Range("B11").Value = 11.15833369610568
$ws.Range("C11").Value = 6.703602501833321
$ws.Range("D11").Value = 6.556990468029647
$ws.Range("E11").Value = 20.69649541077402
$ws.Range("G11").Value = 3.625403660182611
$ws.Range("I11").Value = 20.67119047552665
$ws.Range("K11").Value = 10.7900905558676
$ws.Range("N11").Value = 17.9133754272418
$ws.Range("O11").Value = 22.62010495170834
# Row 12
$ws.Range("B12").Value = 11.24817675141734
$ws.Range("C12").Value = 6.780750940032933
$ws.Range("D12").Value = 6.59406992749351
$ws.Range("E12").Value = 20.95578919455757
$ws.Range("G12").Value = 3.6248772572974
$ws.Range("I12").Value = 20.65881553892313
$ws.Range("K12").Value = 10.85762493735932
$ws.Range("N12").Value = 17.89689652063512
$ws.Range("O12").Value = 22.62017090001595
# Row 13
$ws.Range("B13").Value = 11.22887868352324
$ws.Range("C13").Value = 6.764214692242068
$ws.Range("D13").Value = 6.586091395964607
$ws.Range("E13").Value = 20.90021182098051
$ws.Range("G13").Value = 3.624990191333345
$ws.Range("I13").Value = 20.66145042876786
$ws.Range("K13").Value = 10.84310474598955
$ws.Range("N13").Value = 17.90043430792746
$ws.Range("O13").Value = 22.62011872461665
# Row 14
$ws.Range("B14").Value = 11.16574812105453
$ws.Range("C14").Value = 6.709985174423282
$ws.Range("D14").Value = 6.560044203442136
$ws.Range("E14").Value = 20.71794766133813
$ws.Range("G14").Value = 3.625360156042196
$ws.Range("I14").Value = 20.67015874413585
$ws.Range("K14").Value = 10.79565761947369
$ws.Range("N14").Value = 17.91201463074035
$ws.Range("O14").Value = 22.62009326806575
# Row 15
$ws.Range("B15").Value = 11.12693005271912
$ws.Range("C15").Value = 6.676536527643758
$ws.Range("D15").Value = 6.544069134083102
$ws.Range("E15").Value = 20.60552562589411
$ws.Range("G15").Value = 3.62558804839467
$ws.Range("I15").Value = 20.67558143560236
$ws.Range("K15").Value = 10.76652408465433
$ws.Range("N15").Value = 17.91914084779379
$ws.Range("O15").Value = 22.62018884123858
# Row 16
$ws.Range("B16").Value = 10.90225654085493
$ws.Range("C16").Value = 6.481306227955192
$ws.Range("D16").Value = 6.452238728743106
$ws.Range("E16").Value = 19.94930089601982
$ws.Range("G16").Value = 3.626913639938925
$ws.Range("I16").Value = 20.70801242456502
$ws.Range("K16").Value = 10.59853464967104
$ws.Range("N16").Value = 17.96048411146697
$ws.Range("O16").Value = 22.62243933738634
# Row 17
$ws.Range("B17").Value = 10.76250832883857
$ws.Range("C17").Value = 6.358370289489897
$ws.Range("D17").Value = 6.395682881018335
$ws.Range("E17").Value = 19.53601195021311
$ws.Range("G17").Value = 3.627744409470302
$ws.Range("I17").Value = 20.72911510147369
$ws.Range("K17").Value = 10.49460723753775
$ws.Range("N17").Value = 17.98629927453722
$ws.Range("O17").Value = 22.6253380561877
# Row 18
$ws.Range("B18").Value = 10.68145089941419
$ws.Range("C18").Value = 6.286496131434742
$ws.Range("D18").Value = 6.363085367934743
$ws.Range("E18").Value = 19.29435021196187
$ws.Range("G18").Value = 3.628228717418808
$ws.Range("I18").Value = 20.74169551156212
$ws.Range("K18").Value = 10.43453272419326
$ws.Range("N18").Value = 18.00131416161299
$ws.Range("O18").Value = 22.62756262489162
# Row 19
$ws.Range("B19").Value = 10.65389248025605
$ws.Range("C19").Value = 6.261960738472835
$ws.Range("D19").Value = 6.352037981932088
$ws.Range("E19").Value = 19.21184846491654
$ws.Range("G19").Value = 3.628393808862826
$ws.Range("I19").Value = 20.74603096436198
$ws.Range("K19").Value = 10.41414345156722
$ws.Range("N19").Value = 18.00642661403693
$ws.Range("O19").Value = 22.62841147959111
# Row 20
$ws.Range("B20").Value = 10.77745555853408
$ws.Range("C20").Value = 6.371577564894014
$ws.Range("D20").Value = 6.401710693461841
$ws.Range("E20").Value = 19.58041571645919
$ws.Range("G20").Value = 3.627655303271611
$ws.Range("I20").Value = 20.72682284459621
$ws.Range("K20").Value = 10.50570185809698
$ws.Range("N20").Value = 17.98353396442957
$ws.Range("O20").Value = 22.62497178922064
# Row 21
$ws.Range("B21").Value = 11.18432223641361
$ws.Range("C21").Value = 6.725961939310878
$ws.Range("D21").Value = 6.567699217907517
$ws.Range("E21").Value = 20.77164555044142
$ws.Range("G21").Value = 3.625251222140981
$ws.Range("I21").Value = 20.66758243204838
$ws.Range("K21").Value = 10.80960885882545
$ws.Range("N21").Value = 17.90860634797118
$ws.Range("O21").Value = 22.62007757602846
# Row 22
$ws.Range("B22").Value = 11.44364491020238
$ws.Range("C22").Value = 6.947210154783163
$ws.Range("D22").Value = 6.675303501507739
$ws.Range("E22").Value = 21.51524466694413
$ws.Range("G22").Value = 3.623737275659634
$ws.Range("I22").Value = 20.632828175569
$ws.Range("K22").Value = 11.00511898206737
$ws.Range("N22").Value = 17.86111206263874
$ws.Range("O22").Value = 22.6218529580928
# Row 23
$ws.Range("B23").Value = 11.30586718619612
$ws.Range("C23").Value = 6.830073307317326
$ws.Range("D23").Value = 6.617966222419068
$ws.Range("E23").Value = 21.12155781006186
$ws.Range("G23").Value = 3.624540075818667
$ws.Range("I23").Value = 20.65101360759187
$ws.Range("K23").Value = 10.90107720271721
$ws.Range("N23").Value = 17.88632611038789
$ws.Range("O23").Value = 22.62044984624725
# Row 24
$ws.Range("B24").Value = 10.77070013454662
$ws.Range("C24").Value = 6.365610283124251
$ws.Range("D24").Value = 6.398985774296471
$ws.Range("E24").Value = 19.56035341174463
$ws.Range("G24").Value = 3.627695567348268
$ws.Range("I24").Value = 20.72785777796176
$ws.Range("K24").Value = 10.50068698802531
$ws.Range("N24").Value = 17.98478362082744
$ws.Range("O24").Value = 22.62513564024555
# Row 25
$ws.Range("B25").Value = 10.1679379515202
$ws.Range("C25").Value = 5.820093588982659
$ws.Range("D25").Value = 6.160358700950209
$ws.Range("E25").Value = 17.73209670724076
$ws.Range("G25").Value = 3.631346915082076
$ws.Range("I25").Value = 20.82765925251361
$ws.Range("K25").Value = 10.05771799338562
$ws.Range("N25").Value = 18.09736456564332
$ws.Range("O25").Value = 22.65148412119424
